# "Displays translation when clicked, then shrinks back to normal."
#
# Column E held a secondary translation for some verbs. Most of those rows
# get folded into column D as "first / second" so a single cell shows both
# meanings; the two "gehen" rows whose secondary sense was basically a
# synonym (uebergehen/"to ignore", untergehen/"to descend") just drop the
# extra column instead of merging it. Column E is then cleared everywhere
# and the sheet shrinks back down to a four-column (A:D) range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> combined column-D text (by row number; row 1 is the header).
# Applied in this order (zugehen first) to match the edit order the
# original author made.
$mergeOrder = @(13, 3, 14, 15, 16)
$merged = @{
    3  = "to tackle / to concern something"   # angehen
    13 = "to approach / to close"             # zugehen
    14 = "to put down / to turn off"          # abstellen
    15 = "to lean on / to attempt"            # anstellen
    16 = "to arrange / to establish"          # aufstellen
}

foreach ($r in $mergeOrder) {
    $ws.Cells.Item($r, 4).Value = $merged[$r]
}

for ($r = 2; $r -le 17; $r++) {
    # Column E (5) is no longer used anywhere - clear it out.
    $ws.Cells.Item($r, 5).ClearContents()
}

# Column D now needs to be wide enough for the combined strings (bestFit
# for the new long "x / y" entries).
$ws.Columns.Item(4).ColumnWidth = 26.75

# Mirror the author's final selection noted in the diff.
$ws.Range("E16").Select()
